# Andrea Miele is now a PO — append a new row (16) to the "2026_po" sheet
# using the same school info as the existing MIPS08000T / Liceo Scientifico
# / A. Volta / Milano / MI / Lombardia rows (e.g. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "MIPS08000T"
$ws.Range("B16").Value = "Liceo Scientifico"
$ws.Range("C16").Value = "A. Volta"
$ws.Range("D16").Value = "Milano"
$ws.Range("E16").Value = "MI"
$ws.Range("F16").Value = "Lombardia"
$ws.Range("G16").Value = "Andrea"
$ws.Range("H16").Value = "Miele"
$ws.Range("I16").Value = "M"

# Data_di_nascita column is formatted as text ("@") in this sheet, so force
# the same number format before writing the date-looking string - otherwise
# Excel would auto-convert it to a date serial number.
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "2008-12-29"

$ws.Range("K16").Value = 4
$ws.Range("L16").Value = "E"

# Leave the selection where Excel would land after typing the last cell of
# the new row (matches the author's recorded cursor position).
[void]$ws.Range("C17").Select()
